# Update Summary sheet (sheet1) rows 2-6, columns A..AU
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Summary")
$ws2 = $wb.Worksheets.Item("Cards_telegram")

$ws1.Cells.Item(2,1).Value = 1369
$ws1.Cells.Item(2,2).Value = 45989.70833333334
$ws1.Cells.Item(2,3).Value = 'Трактор'
$ws1.Cells.Item(2,4).Value = 'ХК Сочи'
$ws1.Cells.Item(2,5).Value = 'Трактор – ХК Сочи'
$ws1.Cells.Item(2,6).Value = 897827
$ws1.Cells.Item(2,7).Value = 'https://text.khl.ru/text/897827.html'
$ws1.Cells.Item(2,8).Value = 5.5
$ws1.Cells.Item(2,9).Value = 1
$ws1.Cells.Item(2,10).Value = 6.5
$ws1.Cells.Item(2,11).Value = 41.132791
$ws1.Cells.Item(2,12).Value = 22.883368
$ws1.Cells.Item(2,13).Value = 64.016159
$ws1.Cells.Item(2,14).Value = 0.779471
$ws1.Cells.Item(2,15).Value = 0.07383199999999999
$ws1.Cells.Item(2,16).Value = 0.09566
$ws1.Cells.Item(2,17).Value = 1.28292136590072
$ws1.Cells.Item(2,18).Value = 13.54426265034132
$ws1.Cells.Item(2,19).Value = 10.45369015262388
$ws1.Cells.Item(2,20).Value = 77.94710000000001
$ws1.Cells.Item(2,21).Value = 7.3832
$ws1.Cells.Item(2,22).Value = 9.565999999999999
$ws1.Cells.Item(2,23).Value = 0.053712
$ws1.Cells.Item(2,24).Value = 0.895251
$ws1.Cells.Item(2,25).Value = 1.117005175084976
$ws1.Cells.Item(2,26).Value = 0.113432
$ws1.Cells.Item(2,27).Value = 0.835531
$ws1.Cells.Item(2,28).Value = 1.196843683836985
$ws1.Cells.Item(2,29).Value = 0.203386
$ws1.Cells.Item(2,30).Value = 0.745577
$ws1.Cells.Item(2,31).Value = 1.341243090921528
$ws1.Cells.Item(2,32).Value = 0.9852300000000001
$ws1.Cells.Item(2,33).Value = 0.01477
$ws1.Cells.Item(2,34).Value = 1.014991423322473
$ws1.Cells.Item(2,35).Value = 0.945893
$ws1.Cells.Item(2,36).Value = 0.054107
$ws1.Cells.Item(2,37).Value = 1.057202030250779
$ws1.Cells.Item(2,38).Value = 0.777299
$ws1.Cells.Item(2,39).Value = 0.222701
$ws1.Cells.Item(2,40).Value = 1.286506222187344
$ws1.Cells.Item(2,41).Value = 0.5423789999999999
$ws1.Cells.Item(2,42).Value = 0.457621
$ws1.Cells.Item(2,43).Value = 1.843729200429958
$ws1.Cells.Item(2,44).Value = 0.90032
$ws1.Cells.Item(2,45).Value = 1.110716189799182
$ws1.Cells.Item(2,46).Value = 0.271571
$ws1.Cells.Item(2,47).Value = 3.682278299229299
$ws1.Cells.Item(3,1).Value = 1369
$ws1.Cells.Item(3,2).Value = 45989.79166666666
$ws1.Cells.Item(3,3).Value = 'Ак Барс'
$ws1.Cells.Item(3,4).Value = 'СКА'
$ws1.Cells.Item(3,5).Value = 'Ак Барс – СКА'
$ws1.Cells.Item(3,6).Value = 897823
$ws1.Cells.Item(3,7).Value = 'https://text.khl.ru/text/897823.html'
$ws1.Cells.Item(3,8).Value = 2.240677
$ws1.Cells.Item(3,9).Value = 3.923077
$ws1.Cells.Item(3,10).Value = 6.163754
$ws1.Cells.Item(3,11).Value = 29.335543
$ws1.Cells.Item(3,12).Value = 34.435929
$ws1.Cells.Item(3,13).Value = 63.771472
$ws1.Cells.Item(3,14).Value = 0.125662
$ws1.Cells.Item(3,15).Value = 0.108911
$ws1.Cells.Item(3,16).Value = 0.760551
$ws1.Cells.Item(3,17).Value = 7.957855198866802
$ws1.Cells.Item(3,18).Value = 9.181809000009183
$ws1.Cells.Item(3,19).Value = 1.314836217426576
$ws1.Cells.Item(3,20).Value = 12.5662
$ws1.Cells.Item(3,21).Value = 10.8911
$ws1.Cells.Item(3,22).Value = 76.0551
$ws1.Cells.Item(3,23).Value = 0.246453
$ws1.Cells.Item(3,24).Value = 0.748671
$ws1.Cells.Item(3,25).Value = 1.335700194077238
$ws1.Cells.Item(3,26).Value = 0.398203
$ws1.Cells.Item(3,27).Value = 0.596921
$ws1.Cells.Item(3,28).Value = 1.675263560839709
$ws1.Cells.Item(3,29).Value = 0.557636
$ws1.Cells.Item(3,30).Value = 0.437488
$ws1.Cells.Item(3,31).Value = 2.285776981311487
$ws1.Cells.Item(3,32).Value = 0.593576
$ws1.Cells.Item(3,33).Value = 0.406424
$ws1.Cells.Item(3,34).Value = 1.684704233324797
$ws1.Cells.Item(3,35).Value = 0.322906
$ws1.Cells.Item(3,36).Value = 0.677094
$ws1.Cells.Item(3,37).Value = 3.096876490371811
$ws1.Cells.Item(3,38).Value = 0.928395
$ws1.Cells.Item(3,39).Value = 0.071605
$ws1.Cells.Item(3,40).Value = 1.077127731192003
$ws1.Cells.Item(3,41).Value = 0.803308
$ws1.Cells.Item(3,42).Value = 0.196692
$ws1.Cells.Item(3,43).Value = 1.244852534768731
$ws1.Cells.Item(3,44).Value = 0.380078
$ws1.Cells.Item(3,45).Value = 2.631038892016901
$ws1.Cells.Item(3,46).Value = 0.937003
$ws1.Cells.Item(3,47).Value = 1.06723244215867
$ws1.Cells.Item(4,1).Value = 1369
$ws1.Cells.Item(4,2).Value = 45989.79166666666
$ws1.Cells.Item(4,3).Value = 'Нефтехимик'
$ws1.Cells.Item(4,4).Value = 'Драконы'
$ws1.Cells.Item(4,5).Value = 'Нефтехимик – Драконы'
$ws1.Cells.Item(4,6).Value = 897826
$ws1.Cells.Item(4,7).Value = 'https://text.khl.ru/text/897826.html'
$ws1.Cells.Item(4,8).Value = 1.591473
$ws1.Cells.Item(4,9).Value = 1.848538
$ws1.Cells.Item(4,10).Value = 3.440011
$ws1.Cells.Item(4,11).Value = 27.688566
$ws1.Cells.Item(4,12).Value = 27.29634
$ws1.Cells.Item(4,13).Value = 54.984906
$ws1.Cells.Item(4,14).Value = 0.436582
$ws1.Cells.Item(4,15).Value = 0.157276
$ws1.Cells.Item(4,16).Value = 0.404708
$ws1.Cells.Item(4,17).Value = 2.290520452057116
$ws1.Cells.Item(4,18).Value = 6.358249192502353
$ws1.Cells.Item(4,19).Value = 2.470917303339692
$ws1.Cells.Item(4,20).Value = 43.6582
$ws1.Cells.Item(4,21).Value = 15.7276
$ws1.Cells.Item(4,22).Value = 40.4708
$ws1.Cells.Item(4,23).Value = 0.203053
$ws1.Cells.Item(4,24).Value = 0.7955140000000001
$ws1.Cells.Item(4,25).Value = 1.257048901716374
$ws1.Cells.Item(4,26).Value = 0.341846
$ws1.Cells.Item(4,27).Value = 0.65672
$ws1.Cells.Item(4,28).Value = 1.522718966987453
$ws1.Cells.Item(4,29).Value = 0.496633
$ws1.Cells.Item(4,30).Value = 0.501934
$ws1.Cells.Item(4,31).Value = 1.992293807552387
$ws1.Cells.Item(4,32).Value = 0.853008
$ws1.Cells.Item(4,33).Value = 0.146992
$ws1.Cells.Item(4,34).Value = 1.172321947742577
$ws1.Cells.Item(4,35).Value = 0.660005
$ws1.Cells.Item(4,36).Value = 0.339995
$ws1.Cells.Item(4,37).Value = 1.515140036817903
$ws1.Cells.Item(4,38).Value = 0.84051
$ws1.Cells.Item(4,39).Value = 0.15949
$ws1.Cells.Item(4,40).Value = 1.189753839930518
$ws1.Cells.Item(4,41).Value = 0.639103
$ws1.Cells.Item(4,42).Value = 0.360897
$ws1.Cells.Item(4,43).Value = 1.564693015053912
$ws1.Cells.Item(4,44).Value = 0.736572
$ws1.Cells.Item(4,45).Value = 1.357640529371195
$ws1.Cells.Item(4,46).Value = 0.709293
$ws1.Cells.Item(4,47).Value = 1.409854601694927
$ws1.Cells.Item(5,1).Value = 1369
$ws1.Cells.Item(5,2).Value = 45989.79166666666
$ws1.Cells.Item(5,3).Value = 'Северсталь'
$ws1.Cells.Item(5,4).Value = 'Локомотив'
$ws1.Cells.Item(5,5).Value = 'Северсталь – Локомотив'
$ws1.Cells.Item(5,6).Value = 897825
$ws1.Cells.Item(5,7).Value = 'https://text.khl.ru/text/897825.html'
$ws1.Cells.Item(5,8).Value = 1.5
$ws1.Cells.Item(5,9).Value = 1.40625
$ws1.Cells.Item(5,10).Value = 2.90625
$ws1.Cells.Item(5,11).Value = 22.31792
$ws1.Cells.Item(5,12).Value = 22.352212
$ws1.Cells.Item(5,13).Value = 44.670133
$ws1.Cells.Item(5,14).Value = 0.371378
$ws1.Cells.Item(5,15).Value = 0.261117
$ws1.Cells.Item(5,16).Value = 0.367505
$ws1.Cells.Item(5,17).Value = 2.692674310271475
$ws1.Cells.Item(5,18).Value = 3.829700862065664
$ws1.Cells.Item(5,19).Value = 2.721051414266472
$ws1.Cells.Item(5,20).Value = 37.1378
$ws1.Cells.Item(5,21).Value = 26.1117
$ws1.Cells.Item(5,22).Value = 36.7505
$ws1.Cells.Item(5,23).Value = 0.87056
$ws1.Cells.Item(5,24).Value = 0.129439
$ws1.Cells.Item(5,25).Value = 7.725646829780824
$ws1.Cells.Item(5,26).Value = 0.947361
$ws1.Cells.Item(5,27).Value = 0.052639
$ws1.Cells.Item(5,28).Value = 18.99732137768575
$ws1.Cells.Item(5,29).Value = 0.981255
$ws1.Cells.Item(5,30).Value = 0.018745
$ws1.Cells.Item(5,31).Value = 53.34755934915977
$ws1.Cells.Item(5,32).Value = 0.383108
$ws1.Cells.Item(5,33).Value = 0.616892
$ws1.Cells.Item(5,34).Value = 2.610230013468787
$ws1.Cells.Item(5,35).Value = 0.149413
$ws1.Cells.Item(5,36).Value = 0.850587
$ws1.Cells.Item(5,37).Value = 6.692858051173593
$ws1.Cells.Item(5,38).Value = 0.380173
$ws1.Cells.Item(5,39).Value = 0.619827
$ws1.Cells.Item(5,40).Value = 2.630381431611398
$ws1.Cells.Item(5,41).Value = 0.147469
$ws1.Cells.Item(5,42).Value = 0.852531
$ws1.Cells.Item(5,43).Value = 6.781086194386617
$ws1.Cells.Item(5,44).Value = 0.835571
$ws1.Cells.Item(5,45).Value = 1.196786389187753
$ws1.Cells.Item(5,46).Value = 0.832979
$ws1.Cells.Item(5,47).Value = 1.200510457046336
$ws1.Cells.Item(6,1).Value = 1369
$ws1.Cells.Item(6,2).Value = 45989.79166666666
$ws1.Cells.Item(6,3).Value = 'Торпедо'
$ws1.Cells.Item(6,4).Value = 'Динамо Мн'
$ws1.Cells.Item(6,5).Value = 'Торпедо – Динамо Мн'
$ws1.Cells.Item(6,6).Value = 897824
$ws1.Cells.Item(6,7).Value = 'https://text.khl.ru/text/897824.html'
$ws1.Cells.Item(6,8).Value = 2.383937
$ws1.Cells.Item(6,9).Value = 4.636364
$ws1.Cells.Item(6,10).Value = 7.0203
$ws1.Cells.Item(6,11).Value = 27.808605
$ws1.Cells.Item(6,12).Value = 39.608055
$ws1.Cells.Item(6,13).Value = 67.41665999999999
$ws1.Cells.Item(6,14).Value = 0.154026
$ws1.Cells.Item(6,15).Value = 0.127631
$ws1.Cells.Item(6,16).Value = 0.716317
$ws1.Cells.Item(6,17).Value = 6.492410372274811
$ws1.Cells.Item(6,18).Value = 7.835087086992973
$ws1.Cells.Item(6,19).Value = 1.396029969971395
$ws1.Cells.Item(6,20).Value = 15.4026
$ws1.Cells.Item(6,21).Value = 12.7631
$ws1.Cells.Item(6,22).Value = 71.6317
$ws1.Cells.Item(6,23).Value = 0.318863
$ws1.Cells.Item(6,24).Value = 0.679112
$ws1.Cells.Item(6,25).Value = 1.472511161634605
$ws1.Cells.Item(6,26).Value = 0.485416
$ws1.Cells.Item(6,27).Value = 0.512559
$ws1.Cells.Item(6,28).Value = 1.95099490985428
$ws1.Cells.Item(6,29).Value = 0.64523
$ws1.Cells.Item(6,30).Value = 0.352744
$ws1.Cells.Item(6,31).Value = 2.834917107023791
$ws1.Cells.Item(6,32).Value = 0.575909
$ws1.Cells.Item(6,33).Value = 0.424091
$ws1.Cells.Item(6,34).Value = 1.736385435893518
$ws1.Cells.Item(6,35).Value = 0.305537
$ws1.Cells.Item(6,36).Value = 0.6944630000000001
$ws1.Cells.Item(6,37).Value = 3.272926028598828
$ws1.Cells.Item(6,38).Value = 0.894557
$ws1.Cells.Item(6,39).Value = 0.105443
$ws1.Cells.Item(6,40).Value = 1.117871751045489
$ws1.Cells.Item(6,41).Value = 0.734795
$ws1.Cells.Item(6,42).Value = 0.265205
$ws1.Cells.Item(6,43).Value = 1.360923795072095
$ws1.Cells.Item(6,44).Value = 0.443622
$ws1.Cells.Item(6,45).Value = 2.254171344072205
$ws1.Cells.Item(6,46).Value = 0.925895
$ws1.Cells.Item(6,47).Value = 1.080036073204845

# Ensure new row 6 date cell keeps the same date/time number format as rows 2-5
$ws1.Cells.Item(6,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Update Cards_telegram sheet (sheet2) rows 2-6, columns A..C
$ws2.Cells.Item(2,1).Value = 45989.70833333334
$ws2.Cells.Item(2,2).Value = 'Трактор – ХК Сочи'
$card2 = @'
КХЛ • Регулярный чемпионат • 28.11.2025
Трактор – ХК Сочи
Ожидания модели (60’):
• Голы: λ_total ≈ 9.04 (6.19 : 2.85)
• Броски: SOG λ ≈ 64 (41 : 23)
Исход (60’), честные кф:
• П1: 77.9%  (Kмод 1.28)
• Х:  7.4%  (Kмод 13.54)
• П2: 9.6%  (Kмод 10.45)
Тоталы голов:
• ТМ 4.5: 5.4%  (Kмод 18.62)
• ТБ 4.5: 89.5%  (Kмод 1.12)
• ТМ 5.5: 11.3%  (Kмод 8.82)
• ТБ 5.5: 83.6%  (Kмод 1.20)
• ТМ 6.5: 20.3%  (Kмод 4.92)
• ТБ 6.5: 74.6%  (Kмод 1.34)
Индивидуальные тоталы:
• Трактор ИТБ 1.5: 98.5% (Kмод 1.01)
• Трактор ИТБ 2.5: 94.6% (Kмод 1.06)
• ХК Сочи ИТБ 1.5: 77.7% (Kмод 1.29)
• ХК Сочи ИТБ 2.5: 54.2% (Kмод 1.84)
Фора +1.5:
• Трактор +1.5: 90.0% (Kмод 1.11)
• ХК Сочи +1.5: 27.2% (Kмод 3.68)
'@
$ws2.Cells.Item(2,3).Value = $card2
$ws2.Cells.Item(3,1).Value = 45989.79166666666
$ws2.Cells.Item(3,2).Value = 'Ак Барс – СКА'
$card3 = @'
КХЛ • Регулярный чемпионат • 28.11.2025
Ак Барс – СКА
Ожидания модели (60’):
• Голы: λ_total ≈ 6.30 (2.00 : 4.31)
• Броски: SOG λ ≈ 64 (29 : 34)
Исход (60’), честные кф:
• П1: 12.6%  (Kмод 7.96)
• Х:  10.9%  (Kмод 9.18)
• П2: 76.1%  (Kмод 1.31)
Тоталы голов:
• ТМ 4.5: 24.6%  (Kмод 4.06)
• ТБ 4.5: 74.9%  (Kмод 1.34)
• ТМ 5.5: 39.8%  (Kмод 2.51)
• ТБ 5.5: 59.7%  (Kмод 1.68)
• ТМ 6.5: 55.8%  (Kмод 1.79)
• ТБ 6.5: 43.7%  (Kмод 2.29)
Индивидуальные тоталы:
• Ак Барс ИТБ 1.5: 59.4% (Kмод 1.68)
• Ак Барс ИТБ 2.5: 32.3% (Kмод 3.10)
• СКА ИТБ 1.5: 92.8% (Kмод 1.08)
• СКА ИТБ 2.5: 80.3% (Kмод 1.24)
Фора +1.5:
• Ак Барс +1.5: 38.0% (Kмод 2.63)
• СКА +1.5: 93.7% (Kмод 1.07)
'@
$ws2.Cells.Item(3,3).Value = $card3
$ws2.Cells.Item(4,1).Value = 45989.79166666666
$ws2.Cells.Item(4,2).Value = 'Нефтехимик – Драконы'
$card4 = @'
КХЛ • Регулярный чемпионат • 28.11.2025
Нефтехимик – Драконы
Ожидания модели (60’):
• Голы: λ_total ≈ 6.69 (3.40 : 3.29)
• Броски: SOG λ ≈ 55 (28 : 27)
Исход (60’), честные кф:
• П1: 43.7%  (Kмод 2.29)
• Х:  15.7%  (Kмод 6.36)
• П2: 40.5%  (Kмод 2.47)
Тоталы голов:
• ТМ 4.5: 20.3%  (Kмод 4.92)
• ТБ 4.5: 79.6%  (Kмод 1.26)
• ТМ 5.5: 34.2%  (Kмод 2.93)
• ТБ 5.5: 65.7%  (Kмод 1.52)
• ТМ 6.5: 49.7%  (Kмод 2.01)
• ТБ 6.5: 50.2%  (Kмод 1.99)
Индивидуальные тоталы:
• Нефтехимик ИТБ 1.5: 85.3% (Kмод 1.17)
• Нефтехимик ИТБ 2.5: 66.0% (Kмод 1.52)
• Драконы ИТБ 1.5: 84.1% (Kмод 1.19)
• Драконы ИТБ 2.5: 63.9% (Kмод 1.56)
Фора +1.5:
• Нефтехимик +1.5: 73.7% (Kмод 1.36)
• Драконы +1.5: 70.9% (Kмод 1.41)
'@
$ws2.Cells.Item(4,3).Value = $card4
$ws2.Cells.Item(5,1).Value = 45989.79166666666
$ws2.Cells.Item(5,2).Value = 'Северсталь – Локомотив'
$card5 = @'
КХЛ • Регулярный чемпионат • 28.11.2025
Северсталь – Локомотив
Ожидания модели (60’):
• Голы: λ_total ≈ 2.65 (1.33 : 1.32)
• Броски: SOG λ ≈ 45 (22 : 22)
Исход (60’), честные кф:
• П1: 37.1%  (Kмод 2.69)
• Х:  26.1%  (Kмод 3.83)
• П2: 36.8%  (Kмод 2.72)
Тоталы голов:
• ТМ 4.5: 87.1%  (Kмод 1.15)
• ТБ 4.5: 12.9%  (Kмод 7.73)
• ТМ 5.5: 94.7%  (Kмод 1.06)
• ТБ 5.5: 5.3%  (Kмод 19.00)
• ТМ 6.5: 98.1%  (Kмод 1.02)
• ТБ 6.5: 1.9%  (Kмод 53.35)
Индивидуальные тоталы:
• Северсталь ИТБ 1.5: 38.3% (Kмод 2.61)
• Северсталь ИТБ 2.5: 14.9% (Kмод 6.69)
• Локомотив ИТБ 1.5: 38.0% (Kмод 2.63)
• Локомотив ИТБ 2.5: 14.7% (Kмод 6.78)
Фора +1.5:
• Северсталь +1.5: 83.6% (Kмод 1.20)
• Локомотив +1.5: 83.3% (Kмод 1.20)
'@
$ws2.Cells.Item(5,3).Value = $card5
$ws2.Cells.Item(6,1).Value = 45989.79166666666
$ws2.Cells.Item(6,2).Value = 'Торпедо – Динамо Мн'
$card6 = @'
КХЛ • Регулярный чемпионат • 28.11.2025
Торпедо – Динамо Мн
Ожидания модели (60’):
• Голы: λ_total ≈ 5.76 (1.93 : 3.82)
• Броски: SOG λ ≈ 67 (28 : 40)
Исход (60’), честные кф:
• П1: 15.4%  (Kмод 6.49)
• Х:  12.8%  (Kмод 7.84)
• П2: 71.6%  (Kмод 1.40)
Тоталы голов:
• ТМ 4.5: 31.9%  (Kмод 3.14)
• ТБ 4.5: 67.9%  (Kмод 1.47)
• ТМ 5.5: 48.5%  (Kмод 2.06)
• ТБ 5.5: 51.3%  (Kмод 1.95)
• ТМ 6.5: 64.5%  (Kмод 1.55)
• ТБ 6.5: 35.3%  (Kмод 2.83)
Индивидуальные тоталы:
• Торпедо ИТБ 1.5: 57.6% (Kмод 1.74)
• Торпедо ИТБ 2.5: 30.6% (Kмод 3.27)
• Динамо Мн ИТБ 1.5: 89.5% (Kмод 1.12)
• Динамо Мн ИТБ 2.5: 73.5% (Kмод 1.36)
Фора +1.5:
• Торпедо +1.5: 44.4% (Kмод 2.25)
• Динамо Мн +1.5: 92.6% (Kмод 1.08)
'@
$ws2.Cells.Item(6,3).Value = $card6

# Ensure new row 6 date cell keeps the same date/time number format as rows 2-5
$ws2.Cells.Item(6,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
